# Fix merge conflict in example_sblar.xlsx: add a new "pricing_interest_rate_type"
# column between the existing "pricing_prepenalty_allowed" / "pricing_prepenalty_exists"
# columns on the "invalid" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# --- 1. Insert a new column at Q (shifts old Q->R and old R->S) ---
$ws.Columns.Item(17).Insert()

# --- 2. New header text for the inserted column ---
$ws.Range("Q1").Value = "pricing_interest_rate_type"

# --- 3. New data for the inserted column (rows 2-11) ---
$newQValues = @{2=1; 3=2; 4=3; 5=4; 6=5; 7=6; 8=999; 9=0; 10=10; 11=1000}
foreach ($r in 2..11) {
    $ws.Cells.Item($r, 17).Value = $newQValues[$r]
}

# --- 4. Drop the inherited wrap-text style on the three pricing columns
#        (Q, R, S) for header + data rows, matching the source file which
#        leaves these cells on the default "Normal" style ---
$ws.Range("Q1:S11").Style = "Normal"

# --- 5. Column widths for the new/resized columns ---
$ws.Columns.Item(17).ColumnWidth = 21.5   # pricing_interest_rate_type
$ws.Columns.Item(18).ColumnWidth = 13.83  # pricing_prepenalty_allowed (was Q)

# --- 6. Row 1 is shorter now that the narrow pricing columns no longer wrap ---
$ws.Rows.Item(1).RowHeight = 34

# --- 7. Update the view: scroll right a bit and select R1:S11 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("R1:S11").Select()
